$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Updated "map" (column D) data points - values re-measured/fixed for the task.
$ws.Range("D3").Value  = 43.330399999999997
$ws.Range("D4").Value  = 90.049499999999995
$ws.Range("D5").Value  = 136.0608
$ws.Range("D6").Value  = 184.1953
$ws.Range("D7").Value  = 235.64750000000001
$ws.Range("D8").Value  = 285.97190000000001
$ws.Range("D9").Value  = 335.30459999999999
$ws.Range("D10").Value = 383.70089999999999
$ws.Range("D11").Value = 438.03129999999999
$ws.Range("D12").Value = 489
$ws.Range("D13").Value = 545.14980000000003
$ws.Range("D14").Value = 591.87840000000006
$ws.Range("D15").Value = 646.94870000000003
$ws.Range("D16").Value = 698.91499999999996
$ws.Range("D17").Value = 752.89819999999997
$ws.Range("D18").Value = 804.10159999999996
$ws.Range("D19").Value = 863.16390000000001
$ws.Range("D20").Value = 910.42470000000003
$ws.Range("D21").Value = 966.05039999999997
$ws.Range("D22").Value = 1017.2333

# Restore the view: keep the sheet active and move the selection/scroll
# position from S27 back up to S8 (this also drops the stale topLeftCell
# that was anchoring the view at row 22).
$ws.Activate() | Out-Null
$ws.Range("S8").Select() | Out-Null
